$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two category labels that changed text (this also updates the
# shared-string table: the two now-unused strings "Umgebung_mit_Technik" and
# "Zentrum_Raumschiff" drop out, and the two new strings are appended in the
# order they're first written to a cell).
$ws.Range("C31").Value = "Technische_Bauten"
$ws.Range("C29").Value = "Zentrum_vom_Raumschiff"

# Update the view: zoom to 175%, scroll so row 7 is at the top, and move the
# selection to C33.
$excel.ActiveWindow.Zoom = 175
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C33").Select()
